$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.218464374542236
$ws.Range("B1").Value = 2.236479043960571
$ws.Range("C1").Value = 6.126449584960938
$ws.Range("D1").Value = 1.998753309249878
$ws.Range("E1").Value = 1.162059426307678
